# Auto-generated script to apply market-data value updates described in the commit diff.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ columns (H, I, J, K, L, M, N) on several rows across all 8 sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2248.7
$ws.Range("J17").Value = 2248.7
$ws.Range("L17").Value = 6746.099999999999
$ws.Range("N17").Value = -7082.099999999999
$ws.Range("H18").Value = 111116450
$ws.Range("J18").Value = 500003140
$ws.Range("L18").Value = 500003140
$ws.Range("N18").Value = -500003708
$ws.Range("H38").Value = 2307.8572
$ws.Range("J38").Value = 7837.25
$ws.Range("L38").Value = 23511.75
$ws.Range("N38").Value = -24255.75
$ws.Range("H61").Value = 150
$ws.Range("I61").Value = 150
$ws.Range("K61").Value = 450
$ws.Range("M61").Value = -278
$ws.Range("H116").Value = 998212.5600000001
$ws.Range("I116").Value = 1164015
$ws.Range("J116").Value = 3398
$ws.Range("K116").Value = 1164015
$ws.Range("L116").Value = 3398
$ws.Range("M116").Value = -1160573
$ws.Range("N116").Value = -10282
$ws.Range("H135").Value = 755.1177
$ws.Range("I135").Value = 362.7
$ws.Range("K135").Value = 3264.3
$ws.Range("M135").Value = -729.2999999999997
$ws.Range("H137").Value = 5196.2593
$ws.Range("I137").Value = 1514.0526
$ws.Range("J137").Value = 13941.5
$ws.Range("K137").Value = 4542.1578
$ws.Range("L137").Value = 41824.5
$ws.Range("M137").Value = -1992.1578
$ws.Range("N137").Value = -46924.5
$ws.Range("H138").Value = 346468.8
$ws.Range("I138").Value = 4064.923
$ws.Range("J138").Value = 535883.75
$ws.Range("K138").Value = 12194.769
$ws.Range("L138").Value = 1607651.25
$ws.Range("M138").Value = -7054.769
$ws.Range("N138").Value = -1617931.25
$ws.Range("H141").Value = 4084.4614
$ws.Range("I141").Value = 3538.8
$ws.Range("K141").Value = 10616.4
$ws.Range("M141").Value = -5436.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5639.327
$ws.Range("I32").Value = 4700.9346
$ws.Range("J32").Value = 12833.667
$ws.Range("K32").Value = 4700.9346
$ws.Range("L32").Value = 12833.667
$ws.Range("M32").Value = -4413.9346
$ws.Range("N32").Value = -13407.667
$ws.Range("H63").Value = 3349.625
$ws.Range("I63").Value = 2971
$ws.Range("K63").Value = 2971
$ws.Range("M63").Value = -2285
$ws.Range("H66").Value = 3349.625
$ws.Range("I66").Value = 2971
$ws.Range("K66").Value = 14855
$ws.Range("M66").Value = -11423
$ws.Range("H74").Value = 406837.84
$ws.Range("I74").Value = 1853973.4
$ws.Range("J74").Value = 12164.546
$ws.Range("K74").Value = 1853973.4
$ws.Range("L74").Value = 12164.546
$ws.Range("M74").Value = -1853099.4
$ws.Range("N74").Value = -13912.546
$ws.Range("H77").Value = 406837.84
$ws.Range("I77").Value = 1853973.4
$ws.Range("J77").Value = 12164.546
$ws.Range("K77").Value = 9269867
$ws.Range("L77").Value = 60822.73
$ws.Range("M77").Value = -9265499
$ws.Range("N77").Value = -69558.73000000001
$ws.Range("H122").Value = 3087.869
$ws.Range("I122").Value = 2745.074
$ws.Range("J122").Value = 5732.2856
$ws.Range("K122").Value = 8235.222
$ws.Range("L122").Value = 17196.8568
$ws.Range("M122").Value = -5785.222
$ws.Range("N122").Value = -22096.8568
$ws.Range("H132").Value = 3329.0715
$ws.Range("J132").Value = 3540.9
$ws.Range("L132").Value = 10622.7
$ws.Range("N132").Value = -15682.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3346045.5
$ws.Range("I107").Value = 4050016.5
$ws.Range("K107").Value = 4050016.5
$ws.Range("M107").Value = -4048096.5
$ws.Range("H134").Value = 3925.5
$ws.Range("I134").Value = 3657.1667
$ws.Range("K134").Value = 10971.5001
$ws.Range("M134").Value = -8436.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3513.8704
$ws.Range("I31").Value = 2848.6562
$ws.Range("J31").Value = 4481.4546
$ws.Range("K31").Value = 2848.6562
$ws.Range("L31").Value = 4481.4546
$ws.Range("M31").Value = -2553.6562
$ws.Range("N31").Value = -5071.4546
$ws.Range("H34").Value = 3513.8704
$ws.Range("I34").Value = 2848.6562
$ws.Range("J34").Value = 4481.4546
$ws.Range("K34").Value = 2848.6562
$ws.Range("L34").Value = 4481.4546
$ws.Range("M34").Value = -2646.6562
$ws.Range("N34").Value = -4885.4546
$ws.Range("H99").Value = 6504
$ws.Range("I99").Value = 6604.8
$ws.Range("K99").Value = 6604.8
$ws.Range("M99").Value = -5106.8
$ws.Range("H126").Value = 6504
$ws.Range("I126").Value = 6604.8
$ws.Range("K126").Value = 19814.4
$ws.Range("M126").Value = -17344.4
$ws.Range("H132").Value = 12503165
$ws.Range("I132").Value = 15627832
$ws.Range("K132").Value = 46883496
$ws.Range("M132").Value = -46880966

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4539.727
$ws.Range("I3").Value = 4177.8
$ws.Range("K3").Value = 12533.4
$ws.Range("M3").Value = -12421.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5488.7554
$ws.Range("I102").Value = 857.75
$ws.Range("K102").Value = 857.75
$ws.Range("M102").Value = 764.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1267.0625
$ws.Range("I61").Value = 1119.6
$ws.Range("J61").Value = 1793.7142
$ws.Range("K61").Value = 1119.6
$ws.Range("L61").Value = 1793.7142
$ws.Range("M61").Value = -917.5999999999999
$ws.Range("N61").Value = -2197.7142
$ws.Range("H113").Value = 1267.0625
$ws.Range("I113").Value = 1119.6
$ws.Range("J113").Value = 1793.7142
$ws.Range("K113").Value = 1119.6
$ws.Range("L113").Value = 1793.7142
$ws.Range("M113").Value = 1050.4
$ws.Range("N113").Value = -6133.7142
$ws.Range("H132").Value = 4572.25
$ws.Range("J132").Value = 10336.429
$ws.Range("L132").Value = 31009.287
$ws.Range("N132").Value = -36069.287
$ws.Range("H136").Value = 5114.6113
$ws.Range("I136").Value = 5827.7
$ws.Range("K136").Value = 17483.1
$ws.Range("M136").Value = -14933.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6979.8
$ws.Range("I81").Value = 12499.5
$ws.Range("J81").Value = 5599.875
$ws.Range("K81").Value = 24999
$ws.Range("L81").Value = 11199.75
$ws.Range("M81").Value = -23938
$ws.Range("N81").Value = -13321.75
$ws.Range("H84").Value = 6979.8
$ws.Range("I84").Value = 12499.5
$ws.Range("J84").Value = 5599.875
$ws.Range("K84").Value = 124995
$ws.Range("L84").Value = 55998.75
$ws.Range("M84").Value = -119691
$ws.Range("N84").Value = -66606.75
$ws.Range("H113").Value = 535.2308
$ws.Range("J113").Value = 783.6667
$ws.Range("L113").Value = 2351.0001
$ws.Range("N113").Value = -6691.0001
$ws.Range("H122").Value = 12502831
$ws.Range("I122").Value = 2548.4285
$ws.Range("J122").Value = 41670156
$ws.Range("K122").Value = 7645.2855
$ws.Range("L122").Value = 125010468
$ws.Range("M122").Value = -5195.2855
$ws.Range("N122").Value = -125015368
